$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B13").Value = [double]"-3.542003062378954e-14"
$ws.Range("B14").Value = -0.02090286107912208
$ws.Range("B15").Value = -0.02090286107912208
$ws.Range("B16").Value = 0.0009501300490622643
$ws.Range("B17").Value = 0.0009501300490622644
$ws.Range("B18").Value = -0.02185299112814893
$ws.Range("B19").Value = -1
$ws.Range("B23").Value = 0.9851779712348125
$ws.Range("B25").Value = 0.9642751101557259
$ws.Range("B26").Value = 0.9642751101557259
$ws.Range("B27").Value = 0.9642751101557259
$ws.Range("B28").Value = 1.816624900533283
$ws.Range("B29").Value = -1.816624900533284
$ws.Range("B30").Value = -1.741944678677917
$ws.Range("B31").Value = 1.741944678677917
$ws.Range("B32").Value = 0.5053622964643509
$ws.Range("B40").Value = 0.01653226285347811
$ws.Range("B41").Value = -0.5902445397213584
$ws.Range("B42").Value = 0.4653499447736746
$ws.Range("B43").Value = 0.5902445397213939
$ws.Range("B45").Value = 1.27961139680996
$ws.Range("B47").Value = 0.1494079502131859
$ws.Range("B48").Value = 0.1494079502131858
$ws.Range("B49").Value = 0.1494079502131858
$ws.Range("B50").Value = 0.1494079502131859
$ws.Range("B55").Value = 0.1494079502131859
$ws.Range("B56").Value = 0.1494079502131859
$ws.Range("B57").Value = 0.1494079502131859
$ws.Range("B58").Value = 0.1494079502131859
$ws.Range("B63").Value = 0.1494079502131859
$ws.Range("B64").Value = 0.1494079502131859
$ws.Range("B65").Value = 0.1494079502131859
$ws.Range("B66").Value = 0.1494079502131859
$ws.Range("B69").Value = 0.1494079502131859
$ws.Range("B70").Value = 0.1494079502131859
$ws.Range("B104").Value = 0
$ws.Range("B106").Value = 1.244907896768409
$ws.Range("B108").Value = -1.244907896768409
$ws.Range("B109").Value = 1.244907896768409
$ws.Range("B118").Value = 0.4829273506811085
$ws.Range("B119").Value = 0
$ws.Range("B120").Value = -1.400610458556501
$ws.Range("B121").Value = 1.400610458556501
$ws.Range("B122").Value = 1.048587275383301
$ws.Range("B123").Value = -1.385075832254562
$ws.Range("B124").Value = 1.385075832254562
$ws.Range("B127").Value = 0.06669912944334273
$ws.Range("B129").Value = 0.06669912944334273
$ws.Range("B130").Value = 0.02109288708892035
$ws.Range("B133").Value = 0.02109288708892035
$ws.Range("B136").Value = [double]"4.728572121238923e-16"
$ws.Range("B137").Value = 1.105381299065313
$ws.Range("B138").Value = -0.4914072613689011
$ws.Range("B139").Value = 0.4914072613689011
$ws.Range("B140").Value = -0.01197163861803588
$ws.Range("B141").Value = 0.01197163861803588
$ws.Range("B142").Value = 0.07981092412023919
$ws.Range("B143").Value = 0.2702169859499527
$ws.Range("B144").Value = -0.04902671053100444
$ws.Range("B145").Value = 0.04902671053100444
$ws.Range("B147").Value = 0.05700780294302799
$ws.Range("B148").Value = 0.05700780294302798
$ws.Range("B151").Value = 0.01653226285347811
$ws.Range("B153").Value = 0.4653499447736746
$ws.Range("B165").Value = 0.1494079502131859
$ws.Range("B175").Value = 0.4172971175429629
$ws.Range("B188").Value = 2.052280905949007
$ws.Range("B192").Value = -2.052280905949007
$ws.Range("B199").Value = -0.1841352035059804
$ws.Range("B200").Value = 0
$ws.Range("B210").Value = 0.9999999999999998
$ws.Range("B241").Value = 0
$ws.Range("B242").Value = 0
